$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.51280000000001
$ws.Range("B4").Value = 5.244699999999995
$ws.Range("E4").Value = 13.77330000000001

$ws.Range("B5").Value = 4.806799999999999

$ws.Range("A7").Value = -21.4

$ws.Range("B8").Value = 4.991499999999999

$ws.Range("E9").Value = 13.90620000000001

$ws.Range("A16").Value = -21.43190000000002
$ws.Range("B16").Value = 4.811600000000001

$ws.Range("E18").Value = 12.5787
